# "Code update and other menu items"
#
# 1. Rename the "output" sheet to "finaldata".
# 2. Add two new rows to the "settings" sheet: "Upper Limit" / "Lower Limit".
# 3. Update the active sheet / selections to match the saved view state
#    (voltages!C13 selected but not the active tab; settings becomes the
#    active tab with A6 selected).

$wb = $excel.ActiveWorkbook

# --- 1. Rename output -> finaldata -------------------------------------
$wsOutput = $wb.Worksheets.Item("output")
$wsOutput.Name = "finaldata"

# --- 2. Append Upper Limit / Lower Limit rows on settings ---------------
$wsSettings = $wb.Worksheets.Item("settings")
$wsSettings.Range("A4").Value = "Upper Limit"
$wsSettings.Range("A5").Value = "Lower Limit"

# --- 3. View state: selections + active sheet/tab -----------------------
$wsVoltages = $wb.Worksheets.Item("voltages")
$wsVoltages.Activate()
$wsVoltages.Range("C13").Select()

$wsSettings.Activate()
$wsSettings.Range("A6").Select()
